$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (borders/styles) from the last existing data row (24)
# down onto the two new rows (25,26) before filling in content - this
# reproduces the same cell styles (s="6" for A/C/D/E, s="3" for G) that the
# rest of the table uses, without disturbing row 24's own formula/value.
$ws.Range("A24:G24").Copy()
$ws.Range("A25:G26").PasteSpecial(-4122)

# Row 25 - New Caledonia / Baie de l'Anse-Vata
# (filled in the same order the author typed the values, so new shared
# strings land at the same indices as the target workbook)
$ws.Range("F25").Value = "Fzo8jORoQMo"
$ws.Range("B25").Value = "-22.30265962622518, 166.44579881056939"
$ws.Range("E25").Value = "New Caledonia"
$ws.Range("D25").Value = "Nouméa"
$ws.Range("C25").Value = "Baie de l'Anse-Vata"
$ws.Range("A25").Value = "LIVE, SEA, BEACH, TRAFFIC"
$ws.Range("G25").Formula = "=IsYouTubeVideoValid(F25)"

# Row 26 - Panama / Panama Fruit Feeder Cam
$ws.Range("F26").Value = "WtoxxHADnGk"
$ws.Range("E26").Value = "Panama"
$ws.Range("B26").Value = "8.621350446135837, -80.13963890536638"
$ws.Range("D26").Value = "Provincia de Coclé"
$ws.Range("C26").Value = "Panama Fruit Feeder Cam at Canopy Lodge | Cornell Lab"
$ws.Range("A26").Value = "LIVE, BIRD, NATURE"
$ws.Range("G26").Formula = "=IsYouTubeVideoValid(F26)"

# Matches the author's new selection position after the edit.
$ws.Range("A27").Select()
